$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(48, 2).Value = "EnergySwap"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(49, 2).Value = "RenderToken"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(51, 2).Value = "SynthetixNetwork"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"

$c = $ws.Cells.Item(2, 4)
$c.Formula = "=""29.899.05"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(2, 5)
$c.Formula = "=""  -0.30%  """
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(3, 4)
$c.Formula = "=""1.898.66"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(3, 5)
$c.Formula = "=""  +0.10%  """
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(4, 4)
$c.Formula = "=""1.000"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(4, 5)
$c.Formula = "=""  +0.04%  """
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(5, 4)
$c.Formula = "=""0.7980"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(5, 5)
$c.Formula = "=""  -4.95%  """
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(6, 4)
$c.Formula = "=""244.36"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(6, 5)
$c.Formula = "=""  +1.12%  """
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(7, 4)
$c.Formula = "=""1.000"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(7, 5)
$c.Formula = "=""  +0.05%  """
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(8, 4)
$c.Formula = "=""0.3174"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(8, 5)
$c.Formula = "=""  -3.49%  """
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(9, 4)
$c.Formula = "=""25.49"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(9, 5)
$c.Formula = "=""  -4.28%  """
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(10, 4)
$c.Formula = "=""0.07169"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(10, 5)
$c.Formula = "=""  +1.61%  """
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(11, 4)
$c.Formula = "=""0.08119"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(11, 5)
$c.Formula = "=""  +0.53%  """
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(12, 4)
$c.Formula = "=""0.7719"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(12, 5)
$c.Formula = "=""  +1.81%  """
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(13, 4)
$c.Formula = "=""5.557"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(13, 5)
$c.Formula = "=""  +5.65%  """
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(14, 4)
$c.Formula = "=""1.899.90"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(14, 5)
$c.Formula = "=""  +0.21%  """
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(15, 4)
$c.Formula = "=""92.87"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(15, 5)
$c.Formula = "=""  +0.65%  """
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(16, 4)
$c.Formula = "=""6.173"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(16, 5)
$c.Formula = "=""  +4.98%  """
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(17, 4)
$c.Formula = "=""29.894.92"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(17, 5)
$c.Formula = "=""  -0.28%  """
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(18, 4)
$c.Formula = "=""13.96"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(18, 5)
$c.Formula = "=""  -0.96%  """
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(19, 4)
$c.Formula = "=""245.60"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(19, 5)
$c.Formula = "=""  +0.37%  """
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(20, 4)
$c.Formula = "=""0.000007769"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(20, 5)
$c.Formula = "=""  -0.05%  """
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(21, 4)
$c.Formula = "=""8.266"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(21, 5)
$c.Formula = "=""  +18.24%  """
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(22, 4)
$c.Formula = "=""1.000"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(22, 5)
$c.Formula = "=""  +0.19%  """
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(23, 4)
$c.Formula = "=""2.149.68"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(23, 5)
$c.Formula = "=""  +0.23%  """
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(24, 4)
$c.Formula = "=""1.001"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(24, 5)
$c.Formula = "=""  +0.06%  """
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(25, 4)
$c.Formula = "=""0.1681"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(25, 5)
$c.Formula = "=""  -3.35%  """
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(26, 4)
$c.Formula = "=""9.480"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(26, 5)
$c.Formula = "=""  +2.46%  """
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(27, 4)
$c.Formula = "=""164.60"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(27, 5)
$c.Formula = "=""  -0.76%  """
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(28, 4)
$c.Formula = "=""18.76"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(28, 5)
$c.Formula = "=""  -0.67%  """
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(29, 4)
$c.Formula = "=""2.077"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(29, 5)
$c.Formula = "=""  -1.35%  """
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(30, 4)
$c.Formula = "=""1.412"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(30, 5)
$c.Formula = "=""  +3.85%  """
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(31, 4)
$c.Formula = "=""1.548"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(31, 5)
$c.Formula = "=""  +2.05%  """
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(32, 4)
$c.Formula = "=""4.509"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(32, 5)
$c.Formula = "=""  +4.95%  """
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(33, 4)
$c.Formula = "=""0.05645"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(33, 5)
$c.Formula = "=""  -4.53%  """
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(34, 4)
$c.Formula = "=""4.092"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(34, 5)
$c.Formula = "=""  +0.32%  """
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(35, 4)
$c.Formula = "=""1.289"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(35, 5)
$c.Formula = "=""  +1.12%  """
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(36, 4)
$c.Formula = "=""0.7445"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(36, 5)
$c.Formula = "=""  +1.62%  """
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(37, 4)
$c.Formula = "=""1.002"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(37, 5)
$c.Formula = "=""  +0.42%  """
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(38, 4)
$c.Formula = "=""2.635"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(38, 5)
$c.Formula = "=""  -3.17%  """
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(39, 4)
$c.Formula = "=""0.01937"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(39, 5)
$c.Formula = "=""  +0.94%  """
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(40, 4)
$c.Formula = "=""2.786"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(40, 5)
$c.Formula = "=""  +0.47%  """
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(41, 4)
$c.Formula = "=""1.171.18"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(41, 5)
$c.Formula = "=""  +16.14%  """
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(42, 4)
$c.Formula = "=""74.94"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(42, 5)
$c.Formula = "=""  +3.45%  """
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(43, 4)
$c.Formula = "=""0.4435"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(43, 5)
$c.Formula = "=""  -0.30%  """
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(44, 4)
$c.Formula = "=""5.944"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(44, 5)
$c.Formula = "=""  +1.15%  """
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(45, 4)
$c.Formula = "=""0.8558"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(45, 5)
$c.Formula = "=""  +1.56%  """
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(46, 4)
$c.Formula = "=""104.54"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(46, 5)
$c.Formula = "=""  +2.64%  """
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(47, 4)
$c.Formula = "=""1.000"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(47, 5)
$c.Formula = "=""  +0.09%  """
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(48, 4)
$c.Formula = "=""10.09"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(48, 5)
$c.Formula = "=""  +2.26%  """
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(49, 4)
$c.Formula = "=""1.886"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(49, 5)
$c.Formula = "=""  +0.03%  """
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(50, 4)
$c.Formula = "=""7.483"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(50, 5)
$c.Formula = "=""  -1.26%  """
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(51, 4)
$c.Formula = "=""2.959"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Cells.Item(51, 5)
$c.Formula = "=""  +8.98%  """
$c.Copy()
$c.PasteSpecial(-4163)
